$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted into the dataset at row 762,
# pushing every subsequent record (old rows 762:863) down by one row
# (new rows 763:864). EntireRow.Insert mirrors Excel's own "Insert Row"
# behaviour (existing rows shift down, formatting of the row above is
# carried over - which is why column D, a date, keeps its date style).
$ws.Rows.Item(762).EntireRow.Insert()

# Fill in the new record's values.
$ws.Cells.Item(762, 1).Value = 10
$ws.Cells.Item(762, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(762, 3).Value = "La Araucanía"
$ws.Cells.Item(762, 4).Value2 = 45124
$ws.Cells.Item(762, 5).Value = 9
$ws.Cells.Item(762, 6).Value = 100112032
$ws.Cells.Item(762, 7).Value = "Zapallo italiano"
$ws.Cells.Item(762, 8).Value = "Sin especificar"
$ws.Cells.Item(762, 9).Value = "Primera"
$ws.Cells.Item(762, 10).Value = 450
$ws.Cells.Item(762, 11).Value = 17000
$ws.Cells.Item(762, 12).Value = 18000
$ws.Cells.Item(762, 13).Value = 17444
$ws.Cells.Item(762, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(762, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(762, 16).Value = 349
$ws.Cells.Item(762, 17).Value = 50
$ws.Cells.Item(762, 18).Value = "Hortaliza"
